$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column L ("2020") mirrors the existing "2019" column (K), reusing that
# column's cell formatting (font / border / number format) and then
# overwriting the values for the new year.

# L2 - thin header-row filler cell, same format as K2 (no value)
$ws.Range("K2").Copy()
$ws.Range("L2").PasteSpecial(-4122)
$ws.Range("L2").Value = ""

# L3 - year header "2020", same bold format as K3 ("2019")
$ws.Range("K3").Copy()
$ws.Range("L3").PasteSpecial(-4122)
$ws.Range("L3").Value = 2020
$ws.Range("L3").VerticalAlignment = -4107

# L4 - same format as K4
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("L4").Value = 1004
$ws.Range("L4").VerticalAlignment = -4107

# L5 - same format as K5
$ws.Range("K5").Copy()
$ws.Range("L5").PasteSpecial(-4122)
$ws.Range("L5").Value = 8279
$ws.Range("L5").VerticalAlignment = -4107

# L6 - same format as K6
$ws.Range("K6").Copy()
$ws.Range("L6").PasteSpecial(-4122)
$ws.Range("L6").Value = 1752
$ws.Range("L6").VerticalAlignment = -4107

# L7 - same format as K7
$ws.Range("K7").Copy()
$ws.Range("L7").PasteSpecial(-4122)
$ws.Range("L7").Value = 6527
$ws.Range("L7").VerticalAlignment = -4107

# L8 - same format as K8 (#,##0 number format)
$ws.Range("K8").Copy()
$ws.Range("L8").PasteSpecial(-4122)
$ws.Range("L8").Value = 10324
$ws.Range("L8").VerticalAlignment = -4107

# L9 - same format as K9
$ws.Range("K9").Copy()
$ws.Range("L9").PasteSpecial(-4122)
$ws.Range("L9").Value = 4131
$ws.Range("L9").VerticalAlignment = -4107

# L10 - same format as K10 (bottom border)
$ws.Range("K10").Copy()
$ws.Range("L10").PasteSpecial(-4122)
$ws.Range("L10").Value = 6193
$ws.Range("L10").VerticalAlignment = -4107

$excel.CutCopyMode = 0

# Move the active selection, matching the saved cursor position in the
# source file.
$ws.Range("Q11").Select()
